# Updates the "Date" column (B) with new test-run timestamps, as produced
# by a fresh Katalon test execution on Thu Mar 07 2024.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @{
    2  = "Thu Mar 07 23:38:43 EST 2024"
    3  = "Thu Mar 07 23:39:08 EST 2024"
    4  = "Thu Mar 07 23:39:46 EST 2024"
    5  = "Thu Mar 07 23:40:32 EST 2024"
    6  = "Thu Mar 07 23:41:14 EST 2024"
    7  = "Thu Mar 07 23:41:44 EST 2024"
    8  = "Thu Mar 07 23:42:09 EST 2024"
    9  = "Thu Mar 07 23:42:38 EST 2024"
    10 = "Thu Mar 07 23:43:21 EST 2024"
    13 = "Thu Mar 07 23:43:46 EST 2024"
    14 = "Thu Mar 07 23:44:18 EST 2024"
    15 = "Thu Mar 07 23:45:07 EST 2024"
    16 = "Thu Mar 07 23:45:36 EST 2024"
    17 = "Thu Mar 07 23:46:22 EST 2024"
    18 = "Thu Mar 07 23:46:48 EST 2024"
    19 = "Thu Mar 07 23:47:31 EST 2024"
    20 = "Thu Mar 07 23:47:56 EST 2024"
    21 = "Thu Mar 07 23:48:30 EST 2024"
    22 = "Thu Mar 07 23:48:50 EST 2024"
    23 = "Thu Mar 07 23:49:29 EST 2024"
    24 = "Thu Mar 07 23:50:07 EST 2024"
    25 = "Thu Mar 07 23:50:23 EST 2024"
    26 = "Thu Mar 07 23:50:48 EST 2024"
    27 = "Thu Mar 07 23:51:30 EST 2024"
    28 = "Thu Mar 07 23:52:21 EST 2024"
    29 = "Thu Mar 07 23:52:49 EST 2024"
    30 = "Thu Mar 07 23:53:16 EST 2024"
}

foreach ($row in $dates.Keys) {
    $ws.Cells.Item($row, 2).Value = $dates[$row]
}
